$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember how many data rows we have before making any structural changes.
$lastRow = $ws.UsedRange.Rows.Count

# Shift the existing columns A:E one place to the right (-> B:F), carrying
# their values and formatting with them.
$ws.Range("A1").EntireColumn.Insert()

# New column B header.
$ws.Range("B1").Value = "segments"

# Give the new "segments" header the same look (bold, bordered, centered) as
# the other header cells, by copying the format from a neighbouring header.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("B1").PasteSpecial(-4122) | Out-Null

# Populate the new column A with the zero-based segment index and apply the
# same bold/bordered header-style formatting that the segment-name column
# used to have.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("A2:A" + $lastRow).PasteSpecial(-4122) | Out-Null
for ($i = 2; $i -le $lastRow; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 2
}

# The segment-name column (now B) no longer keeps that bold/bordered style.
$ws.Range("B2:B" + $lastRow).ClearFormats()

$excel.CutCopyMode = 0
